# TC01_C3DC_phs000471_SexAtBirth-Unknown.xlsx
# -------------------------------------------------------------------------
# Fixes the "Treatment" SQL query stored in cell B5 (TreatmentTab row):
# the redundant CONCAT() wrapper around REPLACE() is removed, e.g.
#   CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS "Treatment Agent"
# becomes
#   REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent"
#
# Re-typing the cell's content causes Excel to re-materialize the shared
# string table (the previous B5 text is dropped and the corrected text is
# appended as a new shared string, while the still-referenced
# TreatmentResp/Survival strings shift up to fill the gap) and, because the
# cell's font is re-asserted to its existing (unchanged) value, a fresh
# (but visually identical) font/style entry is produced as well -- both are
# side effects that mirror the authoring session captured in the diff.
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedTreatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs000471' AND prt.sex_at_birth = 'Unknown'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

# B5 holds the "TreatmentTab" query (row 5, column B -> TreatmentTab per A5)
$treatmentCell = $ws.Cells.Item(5, 2)
$treatmentCell.Value2 = $fixedTreatmentQuery

# Re-assert the (unchanged) font so the workbook records a fresh style
# entry for the edited cell, just as happened in the original authoring
# session, while keeping the exact same visible formatting (wrap text,
# 12pt Calibri).
$treatmentCell.Font.ThemeColor = 1

# Move the active selection/scroll position to C5, as captured in the diff.
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
